$d = $word.ActiveDocument

# --- Title / Title Char: drop the manual character-spacing / kerning
# overrides that used to be baked into the big display heading font,
# leaving just the theme fonts + size. ---
foreach ($styleName in @("Title", "TitleChar")) {
    $s = $d.Styles($styleName)
    $s.Font.Spacing = 0
    $s.Font.Kerning = 0
}

# --- Author / Date: rebase onto Title (so they inherit its centered,
# keep-with-next heading formatting) and give them an explicit, smaller
# 12pt run size instead of the inherited 28pt Title size. The explicit
# center alignment is no longer needed once the style is based on Title,
# which already centers. ---
$titleStyle = $d.Styles("Title")

foreach ($styleName in @("Author", "Date")) {
    $s = $d.Styles($styleName)
    $s.BaseStyle = $titleStyle
    $s.Font.Size = 12
    $s.Font.SizeBi = 12
}
